$d = $word.ActiveDocument

# --- 1. Replace the block spanning "Muster eingesetzt?" through the "..." paragraph ---
$startPara = -1
$endPara = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*Muster eingesetzt?*") {
        $startPara = $i
    }
    if ($startPara -ge 1 -and $i -ge $startPara -and $endPara -eq -1) {
        if ($t -like "*$([char]8230)*") {
            $endPara = $i
        }
    }
}

$rStart = $d.Paragraphs($startPara).Range.Start
$rEnd = $d.Paragraphs($endPara).Range.End
$r = $d.Range($rStart, $rEnd)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:body><w:p w:rsidR="00E262C8" w:rsidRPr="00E262C8" w:rsidRDefault="00522037" w:rsidP="00522037"><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">    Muster eingesetzt?</w:t></w:r><w:r><w:rPr><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Default"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t>Ein Strategie</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">muster nutzt man um eine Nutzerklasse (Klient) von </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">zur Laufzeit </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t>austauschbaren Algorithmen abzukoppeln. Der Klient greift dabei a</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">uf einen sogenannten Kontext zu, der die vom Klient gewünschten Operationen durchführt, ohne dass der Klient selbst bemerkt, auf welche Weise dies geschieht. Der Kontext greift dabei auf ein Regelwerk zurück, um zur Laufzeit zu entscheiden, welche Strategie er nutzen möchte, um eine Anfrage des Klienten zu erfüllen. Die konkrete Strategie ist ein Objekt einer von einer Basisklasse „Strategie“ abgeleiteten Klasse, sodass der Kontext lediglich die von der Basisklasse bereitgestellte </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t>Algorithmusschnittstelle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> nutzen und keine etwaigen Besonderheiten der konkreten Strategie bedenken muss. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">Eine Einsatzmöglichkeit bietet sich zum Beispiel in einem Kartenspiel, in welchem eine Künstliche Intelligenz „am Zuge“ ist und anhand eines Regelwerks </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve">entsprechend entscheidet, wie ihr Zug vollzogen werden soll, die eingehenden Daten wären zum Beispiel die oben liegende (n) Karte(n) und das eigene Blatt, konkrete Strategien wären vielleicht „Bedienen“ oder andere regelkonforme </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t>Zugarten</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cstheme="minorHAnsi"/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="de-DE"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00AD1F9F" w:rsidRPr="00FF0BBA" w:rsidRDefault="00E262C8" w:rsidP="00CA5284"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xml)

Write-Output "block replaced: startPara=$startPara endPara=$endPara"

# --- 2. Update footer NUMPAGES cached field text from 5 to 4 ---
$footer = $d.Sections(1).Footers(1)
$numPagesField = $footer.Range.Fields(2)
$numPagesField.Result.Text = "4"
Write-Output "footer updated"
